# Apply the commit's changes to the workbook via Excel COM interop.
#
# Summary of the change (from the OOXML diff):
#  - testCitizen!A1:B8: language name / citizen-code labels get new text
#    ("Turkish012"/"trksh12" -> "Turkish - TR"/"TR-CTZ", etc.)
#  - testCitizen!A1:A8: wrap text turned on (new cell style)
#  - testCitizen sheet view: selection moves from B8 to C13
#  - workbook bookView: windowHeight 13980 -> 12180

$wb  = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item("testCitizen")

# --- Relabel the language / citizenship columns (A and B) ---
$ws2.Range("A1").Value = "Turkish - TR"
$ws2.Range("B1").Value = "TR-CTZ"

$ws2.Range("A2").Value = "Kurdish - KU"
$ws2.Range("B2").Value = "KU-CTZ"

$ws2.Range("A3").Value = "German - DE"
$ws2.Range("B3").Value = "DE-CTZ"

$ws2.Range("A4").Value = "English - EN"
$ws2.Range("B4").Value = "EN-CTZ"

$ws2.Range("A5").Value = "Norwegian - NO"
$ws2.Range("B5").Value = "NO-CTZ"

$ws2.Range("A6").Value = "Italian - IT"
$ws2.Range("B6").Value = "IT-CTZ"

$ws2.Range("A7").Value = "Spanish - ES"
$ws2.Range("B7").Value = "ES-CTZ"

$ws2.Range("A8").Value = "American - US"
$ws2.Range("B8").Value = "US-CTZ"

# --- Turn on wrap text for the label column (creates the new cell style
#     that the diff inserts into cellXfs) ---
$ws2.Range("A1:A8").WrapText = $true

# --- Move the active selection to C13 ---
$ws2.Range("C13").Select()

# --- Workbook window view height (bookViews/workbookView@windowHeight) ---
$excel.ActiveWindow.Height = 12180
